$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("H7").Value = "exponent"
$ws.Range("H9").Value = "bits in register"
$ws.Range("H10").Value = "result of shift"
$ws.Range("I1").Value = "nth register mantissa check"
$ws.Range("H11").Value = "more than 22 bits read in?"
$ws.Range("I14").Value = "final mantissa"

$ws.Columns.Item(8).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(9).EntireColumn.AutoFit() | Out-Null

$excel.ActiveWindow.FreezePanes = $false
$ws.Range("G1").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("I16").Select()
